# Atualização da distribuição de Tarefas. (Participação do grupo inteiro)
#
# The "Website" group of tasks (rows 29-42, column J = Hugo Bastos) is being
# updated to reflect that the whole group participates, so Hugo Bastos'
# individual "X" marks are cleared from most rows; the first row of that
# block (Registo, row 29) keeps a (blank/space) marker instead of "X".

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Row 29 (Registo): replace the "X" mark with a single space.
$ws.Range("J29").Value = " "

# Remaining rows in the Website task block: clear Hugo Bastos' "X" mark,
# since the whole group now participates instead of a single person.
$rowsToClear = @(30, 32, 34, 35, 36, 37, 39, 40, 41, 42)
foreach ($r in $rowsToClear) {
    $ws.Range("J$r").ClearContents()
}

# Update the view: scroll position and current selection.
$win = $excel.ActiveWindow
$win.ScrollRow = 6
$win.ScrollColumn = 1
$ws.Range("K31").Select()
